$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Doa"
$ws.Range("B3").Value = "Dob"
$ws.Range("B4").Value = "Doc"
$ws.Range("B5").Value = "Dod"
$ws.Range("B7").Value = "Dof"
$ws.Range("B8").Value = "Dog"
$ws.Range("B9").Value = "Doh"
$ws.Range("B10").Value = "Doi"
$ws.Range("B11").Value = "Doj"

$ws.Range("B12").Select()
